# Update checkout page verification sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Checkout")

# Rows 49, 50, 53, 54 (column C): "Not found this element" -> "Element not found"
$ws.Cells.Item(49, 3).Value = "Element not found"
$ws.Cells.Item(50, 3).Value = "Element not found"
$ws.Cells.Item(53, 3).Value = "Element not found"
$ws.Cells.Item(54, 3).Value = "Element not found"

# Row 72 (column B): clear "Different" status
$ws.Cells.Item(72, 2).Value = ""

# Row 107 (column B): "Different" -> "Same"
$ws.Cells.Item(107, 2).Value = "Same"

# Row 111: Payment_Input| removed from locator list, replaced by City_Input|
$ws.Cells.Item(111, 1).Value = "City_Input|"
$ws.Cells.Item(111, 2).Value = "Same"
$ws.Cells.Item(111, 3).Value = ""

# Rows 112-115, 117: clear the "Elements hidden - Defect raised" note, mark Same
$ws.Cells.Item(112, 2).Value = "Same"
$ws.Cells.Item(112, 3).Value = ""

$ws.Cells.Item(113, 2).Value = "Same"
$ws.Cells.Item(113, 3).Value = ""

$ws.Cells.Item(114, 2).Value = "Same"
$ws.Cells.Item(114, 3).Value = ""

$ws.Cells.Item(115, 2).Value = "Same"
$ws.Cells.Item(115, 3).Value = ""

# Row 116: clear the note, leave column B blank (header row for PayPal section)
$ws.Cells.Item(116, 3).Value = ""

$ws.Cells.Item(117, 2).Value = "Same"
$ws.Cells.Item(117, 3).Value = ""

# Row 118: mark Different, with new note
$ws.Cells.Item(118, 2).Value = "Different"
$ws.Cells.Item(118, 3).Value = "Element locators not matching in application"

# Update frozen pane / selection to match the reviewed area (rows around 71-79)
$ws.Activate()
$ws.Range("A72").Select()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A72").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("B79").Select()
